# Result Data.xlsx — roster refresh + header bolding + print setup
# Applies the changes described by the upstream diff:
#   * B2:B11 student names replaced (old Y/Z/A..G placeholder strings drop out
#     of the shared-string table as a side effect of no longer being referenced)
#   * Various score cells in D:J (rows 3-11) updated to new figures
#   * Header row (row 1) made bold
#   * Page setup: Letter/A4-ish (paperSize 9), portrait orientation

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Roster names (column B) ---
# Inserted in this order so the rebuilt shared-strings table lands in the same
# sequence as the target workbook (row 9 / "Praveen" added last).
$ws.Range("B2").Value  = "Aditya"
$ws.Range("B3").Value  = "Bimal"
$ws.Range("B4").Value  = "Deepak"
$ws.Range("B5").Value  = "Faisal"
$ws.Range("B6").Value  = "Gurveer"
$ws.Range("B7").Value  = "Neena"
$ws.Range("B8").Value  = "Piyush"
$ws.Range("B10").Value = "Sangharsh"
$ws.Range("B11").Value = "Vivek"
$ws.Range("B9").Value  = "Praveen"

# --- Score updates ---
# Row 2 (Aditya) is unchanged apart from the name.

# Row 3 (Bimal)
$ws.Range("E3").Value = 55
$ws.Range("F3").Value = 70
$ws.Range("G3").Value = 45
$ws.Range("H3").Value = 25
$ws.Range("I3").Value = 37
$ws.Range("J3").Value = 581

# Row 4 (Deepak)
$ws.Range("D4").Value = 28
$ws.Range("E4").Value = 55
$ws.Range("G4").Value = 53
$ws.Range("H4").Value = 28
$ws.Range("I4").Value = 35
$ws.Range("J4").Value = 550

# Row 5 (Faisal)
$ws.Range("D5").Value = 41
$ws.Range("E5").Value = 47
$ws.Range("F5").Value = 65
$ws.Range("H5").Value = 23
$ws.Range("I5").Value = 28

# Row 6 (Gurveer)
$ws.Range("G6").Value = 73
$ws.Range("H6").Value = 35
$ws.Range("I6").Value = 30
$ws.Range("J6").Value = 680

# Row 7 (Neena)
$ws.Range("E7").Value = 30
$ws.Range("F7").Value = 41
$ws.Range("G7").Value = 42
$ws.Range("H7").Value = 21
$ws.Range("I7").Value = 30
$ws.Range("J7").Value = 521

# Row 8 (Piyush)
$ws.Range("D8").Value = 65
$ws.Range("F8").Value = 81
$ws.Range("G8").Value = 58
$ws.Range("H8").Value = 33
$ws.Range("I8").Value = 35
$ws.Range("J8").Value = 772

# Row 9 (Praveen)
$ws.Range("D9").Value = 30
$ws.Range("F9").Value = 40
$ws.Range("H9").Value = 22
$ws.Range("I9").Value = 28
$ws.Range("J9").Value = 554

# Row 10 (Sangharsh)
$ws.Range("E10").Value = 30
$ws.Range("H10").Value = 26
$ws.Range("J10").Value = 530

# Row 11 (Vivek)
$ws.Range("G11").Value = 51
$ws.Range("H11").Value = 28
$ws.Range("I11").Value = 29

# --- Header formatting: bold row 1 ---
$ws.Range("A1:J1").Font.Bold = $true

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
